# Applies the "Checklist PvA" edit:
#  - Marks a set of checklist rows in column B with the checkmark value "v"
#    (these rows previously had an empty B cell, style preserved).
#  - Moves the active selection from F55 to D112.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$rows = @(65,66,67,68,71,72,73,74,75,76,77,78,89,90,91,92,93,94,95,98,99,100,103,104,105,106,107,110)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "v"
}

$ws.Activate()
$ws.Range("D112").Select() | Out-Null
